# DRILL-8417: Allow Excel Reader to Ignore Formula Errors
#
# Adds a second worksheet ("Sheet with Errors") after Sheet1 containing a
# field1 / field2 / result table where result = field1/field2. One row
# (field2 = 0) produces a #DIV/0! error, which is the scenario the reader
# needs to tolerate.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet with Errors"

# Header row - write B1 before A1 so the shared-string table order matches
# (field2, field1, result).
$ws2.Range("B1").Value = "field2"
$ws2.Range("A1").Value = "field1"
$ws2.Range("C1").Value = "result"

# Data rows: result = field1 / field2.
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Formula = "=A2/B2"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 3
$ws2.Range("C3").Formula = "=A3/B3"

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 4
$ws2.Range("C4").Formula = "=A4/B4"

# field2 = 0 here -> #DIV/0! error in the result column.
$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = 0
$ws2.Range("C5").Formula = "=A5/B5"

$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = 6
$ws2.Range("C6").Formula = "=A6/B6"

# Match the saved selections / active sheet: Sheet1's cursor moves to C4
# and the new sheet (now the active tab) is left selected at E5.
$ws1.Range("C4").Select() | Out-Null
$ws2.Range("E5").Select() | Out-Null
